$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A; this shifts B:F left into A:E, matching the target layout.
$ws.Range("A:A").Delete()

# Fix the mislabeled header text (now in D1 after the column shift).
$ws.Range("D1").Value = "MODELCONDITION"
